$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.978.48'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '1.704.15'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3993'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4032'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.470'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.77'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.002'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08806'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.469'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.972'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001349'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '1.730.20'
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '96.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07194'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.306'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').Value = '24.991.00'
$ws.Range('E24').Value = '  +2.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.410'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.945'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.104'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.93'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '152.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.418'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.658'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +19.96%  '
$ws.Range('D33').Value = '1.930.68'
$ws.Range('E33').Value = '  +2.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08588'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.03157'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.045'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.207'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2899'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09625'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8254'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.76%  '
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.479'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.683'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7377'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.31%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.09079'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +11.30%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.251'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.400'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.001'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '139.75'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.00%  '
